# Actualización desde MV -datos-
# Update the last existing row (139 -> "01-06-2021") with the new monthly
# figures, then append the new month row (140 -> "01-07-2021").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 139 values ---
$ws.Cells.Item(139, 2).Value = 0.4
$ws.Cells.Item(139, 3).Value = 0.3
$ws.Cells.Item(139, 4).Value = 0.5

# --- Append new row 140 ---
# Column A holds a text label that looks like a date ("01-07-2021").
# Format the cell as Text first so Excel stores it as a literal string
# (shared string) instead of auto-converting it to a date serial number,
# then restore the cell style back to Normal/General.
$a140 = $ws.Cells.Item(140, 1)
$a140.NumberFormat = "@"
$a140.Value = "01-07-2021"
$a140.Style = "Normal"

$ws.Cells.Item(140, 2).Value = 1.3
$ws.Cells.Item(140, 3).Value = 0.5
$ws.Cells.Item(140, 4).Value = 1.4
